$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose values are numeric- or percent-looking strings need NumberFormat
# forced to text ("@") before assignment so Excel keeps them as text (matching
# the source workbook, which stores them as inline strings), then the format is
# restored to "General" to avoid leaving a stray custom number format behind.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.NumberFormat = "General"
}

Set-TextValue "D2" '305.07'
Set-TextValue "E2" '0.26%'
Set-TextValue "D3" '35.72'
Set-TextValue "E3" '-0.28%'
Set-TextValue "D4" '5.024'
Set-TextValue "E4" '-1.22%'
Set-TextValue "D5" '0.08008'
Set-TextValue "E5" '-0.40%'
Set-TextValue "D6" '1.858'
Set-TextValue "E6" '-3.01%'
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
Set-TextValue "D7" '7.782'
Set-TextValue "E7" '0.64%'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue "D8" '0.9201'
Set-TextValue "E8" '-0.94%'
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue "D9" '0.1267'
Set-TextValue "E9" '-7.78%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue "D10" '0.1914'
Set-TextValue "E10" '1.08%'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue "D11" '0.09137'
Set-TextValue "E11" '0.32%'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue "D12" '0.03463'
Set-TextValue "E12" '-3.49%'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue "D13" '0.09848'
Set-TextValue "E13" '0.42%'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue "D14" '0.001405'
Set-TextValue "E14" '-1.87%'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue "D15" '0.006239'
Set-TextValue "E15" '5.49%'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue "D16" '3.847'
Set-TextValue "E16" '8.25%'
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue "D17" '4.144'
Set-TextValue "E17" '-0.93%'
Set-TextValue "E18" '13.49%'
Set-TextValue "E19" '-1.07%'
Set-TextValue "D20" '0.1348'
Set-TextValue "E20" '1.09%'
Set-TextValue "D21" '5.209'
Set-TextValue "E21" '6.25%'
Set-TextValue "E22" '-8.28%'
Set-TextValue "D23" '0.04425'
Set-TextValue "E23" '-0.81%'
Set-TextValue "D24" '0.001234'
Set-TextValue "E24" '0.75%'
Set-TextValue "D25" '0.004885'
Set-TextValue "E25" '2.10%'
Set-TextValue "E27" '-15.64%'
Set-TextValue "D39" '0.01928'
Set-TextValue "E39" '-1.49%'
Set-TextValue "D40" '0.05210'
Set-TextValue "E40" '6.91%'
Set-TextValue "D41" '0.007616'
Set-TextValue "E41" '-0.12%'
Set-TextValue "D42" '0.01015'
Set-TextValue "E42" '10.57%'
Set-TextValue "D43" '0.1348'
Set-TextValue "E43" '-1.65%'
Set-TextValue "D44" '0.002153'
Set-TextValue "E44" '2.27%'
Set-TextValue "D45" '0.009622'
Set-TextValue "E45" '-15.39%'
Set-TextValue "D46" '0.00006188'
Set-TextValue "E46" '-3.35%'
Set-TextValue "E47" '-0.09%'
Set-TextValue "D48" '64.97'
Set-TextValue "E48" '0.47%'
Set-TextValue "E50" '-0.09%'
Set-TextValue "D51" '0.0002003'
Set-TextValue "E51" '-0.09%'

Write-Host "Applied crypto price/volume update."